# Update cryptos worksheet with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.697.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "'2.096.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'343.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.47%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "'0.5140"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("D8").Value = "'0.4403"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").Value = "'52.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("D10").Value = "'0.09176"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "'24.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").Value = "'2.105.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").Value = "'6.754"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "'8.212"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "'99.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "'0.00001150"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'20.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.47%  "
$ws.Range("D20").Value = "'0.06629"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'6.189"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").Value = "'29.755.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").Value = "'12.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").Value = "'2.320"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").Value = "'2.355.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").Value = "'21.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").Value = "'162.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").Value = "'2.523"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "'132.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.21%  "
$ws.Range("D31").Value = "'1.128"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.29%  "
$ws.Range("D32").Value = "'0.1048"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.51%  "
$ws.Range("D33").Value = "'1.655"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("D34").Value = "'6.158"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("D35").Value = "'3.943"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").Value = "'6.026"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").Value = "'0.02563"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").Value = "'0.06719"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").Value = "'12.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").Value = "'0.2232"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.00%  "
$ws.Range("D42").Value = "'0.6853"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("D43").Value = "'1.287"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "'0.6648"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("D45").Value = "'14.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.62%  "
$ws.Range("D46").Value = "'2.298"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").Value = "'3.607"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.03%  "
$ws.Range("D48").Value = "'0.00000000352"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.02%  "
$ws.Range("D49").Value = "'1.218"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("D50").Value = "'82.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "'0.3278"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.71%  "
